# Record a new entry on the "新题" (New Problems) sheet for problem 102 (bfs),
# marked as done.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("新题")

# New row 20: reuse the date style/format already used by the rows above
# (column A), then fill in the date, the problem title and the "done"
# marker in column F.
$ws.Range("A19").Copy()
$ws.Range("A20").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A20").Value = 43550
$ws.Range("B20").Value = "102 bfs"
$ws.Range("F20").Value = "done"

# Update the active selection to match the author's final cursor position.
$ws.Range("F5").Select()
